$d = $word.ActiveDocument

# Shared run-properties block used by every "Courier New" maze-format run.
$rPr = '<w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:highlight w:val="white"/></w:rPr>'

function Wrap-Body([string]$innerXml) {
    return @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
$innerXml
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

function Set-ParagraphRuns([int]$paraIndex, [string[]]$pieces) {
    # Replaces the run content of a paragraph (but not its paragraph mark /
    # pPr / paragraph attributes) with one <w:r> per entry in $pieces, all
    # sharing the same Courier-New run formatting.
    $p = $d.Paragraphs($paraIndex)
    $full = $p.Range
    $sub = $d.Range($full.Start, $full.End - 1)
    $runs = ""
    foreach ($piece in $pieces) {
        $runs += "<w:r>$rPr<w:t>$piece</w:t></w:r>"
    }
    $xml = Wrap-Body "<w:p>$runs</w:p>"
    $sub.InsertXML($xml)
}

# ---------------------------------------------------------------------
# 1. Drop the stray _GoBack bookmark from the "trap door connection."
#    paragraph (it moves down to the "start:A" paragraph, see below).
# ---------------------------------------------------------------------
$pTrap = $d.Paragraphs(10)
$fullTrap = $pTrap.Range
$xmlTrap = Wrap-Body @"
<w:p w14:paraId="4A013D92" w14:textId="7004E2B0" w:rsidR="000E3C08" w:rsidRDefault="000E3C08" w:rsidP="00175AE1">
<w:r><w:t>A is the name of the room, there is no room set to the north connection, room c is set to the east connection, room f is set to the south connection and there is no room set to the trap door connection.</w:t></w:r>
</w:p>
<w:p/>
"@
$fullTrap.InsertXML($xmlTrap)

# ---------------------------------------------------------------------
# 2. Split each "Room:connections" line into one run per character (with
#    ';' separators) instead of a single run holding the raw string.
# ---------------------------------------------------------------------
Set-ParagraphRuns 32 @("A:-", ";", "C", ";", "F", ";", "-", ";", "-", ";")
Set-ParagraphRuns 33 @("B:-", ";", "-", ";", "C", ";", "-", ";", "-", ";")
Set-ParagraphRuns 34 @("C:B", ";", "D", ";", "G", ";", "A", ";", "H", ";")
Set-ParagraphRuns 35 @("D:-", ";", "E", ";", "-", ";", "C", ";", "-", ";")
Set-ParagraphRuns 36 @("E:-", ";", "-", ";", "-", ";-;", "-", ";")
Set-ParagraphRuns 37 @("F:-", ";", "-", ";", "J", ";", "-", ";", "K", ";")
Set-ParagraphRuns 38 @("G:C", ";", "-", ";", "-", ";", "-", ";", "-", ";")
Set-ParagraphRuns 39 @("H:-", ";", "I", ";", "L", ";", "-", ";", "M", ";")
Set-ParagraphRuns 40 @("I:-", ";", "-", ";", "-", ";", "H", ";", "-", ";")
Set-ParagraphRuns 41 @("J:F", ";", "-", ";", "-", ";", "-", ";", "-", ";")
Set-ParagraphRuns 42 @("K:L", ";", "-", ";", "-", ";", "-", ";", "-", ";")
Set-ParagraphRuns 43 @("L:H", ";", "-", ";", "K", ";", "-", ";", "-", ";")
Set-ParagraphRuns 44 @("M:N", ";", "-", ";", "-", ";", "-", ";", "H", ";")
Set-ParagraphRuns 45 @("N:M", ";", "-", ";", "-", ";", "-", ";", "-", ";")

# ---------------------------------------------------------------------
# 3. "start:A" gains a trailing ';' run and the _GoBack bookmark that was
#    removed from the "trap door connection." paragraph above.
# ---------------------------------------------------------------------
$pStart = $d.Paragraphs(46)
$fullStart = $pStart.Range
$subStart = $d.Range($fullStart.Start, $fullStart.End - 1)
$xmlStart = Wrap-Body @"
<w:p>
<w:r>$rPr<w:t>start:A</w:t></w:r>
<w:r>$rPr<w:t>;</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
"@
$subStart.InsertXML($xmlStart)

# ---------------------------------------------------------------------
# 4. "finish:N" gains a trailing ';' run.
# ---------------------------------------------------------------------
$pFinish = $d.Paragraphs(47)
$fullFinish = $pFinish.Range
$subFinish = $d.Range($fullFinish.Start, $fullFinish.End - 1)
$xmlFinish = Wrap-Body @"
<w:p>
<w:r>$rPr<w:t>finish:N</w:t></w:r>
<w:r>$rPr<w:t>;</w:t></w:r>
</w:p>
"@
$subFinish.InsertXML($xmlFinish)

Write-Output "done"
